$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 'Objectives
id="Par1">To evaluate the prevalence of acute pulmonary embolism (APE) in non-hospitalized COVID-19 patients referred to CT pulmonary angiography (CTPA) by the emergency department.


Methods
id="Par2">From March 14 to April 6, 2020, 72 non-hospitalized patients referred by the emergency department to CTPA for COVID-19 pneumonia were retrospectively identified.

 Relevant clinical and laboratory data and CT scan findings were collected for each patient.

 CTPA scans were reviewed by two radiologists to determinate the presence or absence of APE.

 Clinical classification, lung involvement of COVID-19 pneumonia, and CT total severity score were compared between APE group and non-APE group.


Results
id="Par3">APE was identified in 13 (18%) CTPA scans.

 The mean age and D-dimer of patients from the APE group were higher in comparison with those from the non-APE group (74.4 vs.

 59.6 years, p = 0.008, and 7.29 vs.

 3.29 μg/ml, p = 0.011).

 There was no significant difference between APE and non-APE groups concerning clinical type, COVID-19 pneumonia lung lesions (ground-glass opacity: 85% vs.

 97%; consolidation: 69% vs.

 68%; crazy paving: 38% vs.

 37%; linear reticulation: 69% vs.

 78%), CT severity score (6.3 vs.

 7.1, p = 0.365), quality of CTPA (1.8 vs.

 2.0, p = 0.518), and pleural effusion (38% vs.

 19%, p = 0.146).


Conclusions
id="Par4">Non-hospitalized patients with COVID-19 pneumonia referred to CT scan by the emergency departments are at risk of APE.

 The presence of APE was not limited to severe or critical clinical type of COVID-19 pneumonia.


Key Points
id="Par5">
• Acute pulmonary embolism was found in 18% of non-hospitalized COVID-19 patients referred by the emergency department to CTPA.

 Two (15%) patients had main, four (30%) lobar, and seven (55%) segmental acute pulmonary embolism.



'
$ws.Range("E4").Value = '[Alban%Gervaise%alban.gervaise@hotmail.fr%1,   Caroline%Bouzad%NULL%2,   Caroline%Bouzad%NULL%0,   Evelyne%Peroux%NULL%1,   Carole%Helissey%NULL%1]'
$ws.Range("D9").Value = 'Purpose
id="Par1">Little evidence of increased thrombotic risk is available in COVID-19 patients.

 Our purpose was to assess thrombotic risk in severe forms of SARS-CoV-2 infection.


Methods
id="Par2">All patients referred to 4 intensive care units (ICUs) from two centers of a French tertiary hospital for acute respiratory distress syndrome (ARDS) due to COVID-19 between March 3rd and 31st 2020 were included.

 Medical history, symptoms, biological data and imaging were prospectively collected.

 Propensity score matching was performed to analyze the occurrence of thromboembolic events between non-COVID-19 ARDS and COVID-19 ARDS patients.


Results
id="Par3">150 COVID-19 patients were included (122 men, median age 63 [53; 71] years, SAPSII 49 [37; 64] points).

 Sixty-four clinically relevant thrombotic complications were diagnosed in 150 patients, mainly pulmonary embolisms (16.7%).

 28/29 patients (96.6%) receiving continuous renal replacement therapy experienced circuit clotting.

 Three thrombotic occlusions (in 2 patients) of centrifugal pump occurred in 12 patients (8%) supported by ECMO.

 Most patients (&gt; 95%) had elevated D-dimer and fibrinogen.

 No patient developed disseminated intravascular coagulation.

 Von Willebrand (vWF) activity, vWF antigen and FVIII were considerably increased, and 50/57 tested patients (87.7%) had positive lupus anticoagulant.

 Comparison with non-COVID-19 ARDS patients (n = 145) confirmed that COVID-19 ARDS patients (n = 77) developed significantly more thrombotic complications, mainly pulmonary embolisms (11.7 vs.

 2.1%, p &lt; 0.008).

 Coagulation parameters significantly differed between the two groups.


Conclusion
id="Par4">Despite anticoagulation, a high number of patients with ARDS secondary to COVID-19 developed life-threatening thrombotic complications.

 Higher anticoagulation targets than in usual critically ill patients should therefore probably be suggested.


Electronic supplementary material
The online version of this article (10.1007/s00134-020-06062-x) contains supplementary material, which is available to authorized users.


'
$ws.Range("E9").Value = '[Julie%Helms%NULL%1,   Charles%Tacquard%NULL%1,   François%Severac%NULL%0,   Ian%Leonard-Lorant%NULL%1,   Mickaël%Ohana%NULL%0,   Xavier%Delabranche%NULL%0,   Hamid%Merdji%NULL%1,   Raphaël%Clere-Jehl%NULL%1,   Malika%Schenck%NULL%1,   Florence%Fagot Gandet%NULL%1,   Samira%Fafi-Kremer%NULL%0,   Vincent%Castelain%NULL%1,   Francis%Schneider%NULL%1,   Lélia%Grunebaum%NULL%1,   Eduardo%Anglés-Cano%NULL%1,   Laurent%Sattler%NULL%1,   Paul-Michel%Mertes%NULL%1,   Ferhat%Meziani%Ferhat.Meziani@chru-strasbourg.fr%1,   NULL%NULL%NULL%0]'
